# Added get_form_type in the import_utils
# Adds a new "Form Tag" column (S) to the CapitalCommitment template sheet,
# populating the header and a default value for the two existing data rows,
# and updates the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Form Tag" column header and default values (column S)
$ws.Range("S1").Value = "Form Tag"
$ws.Range("S2").Value = "Default"
$ws.Range("S3").Value = "Default"

# Update the sheet's active selection to the new column
$ws.Range("S4:S5").Select() | Out-Null
